# Daily attendance processing - 2026-01-11 18:39:35
# Swap the order of names in the "Recorded By" column (G) from
# "System, <email>" to "<email>, System" for every matching row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "System, *") {
        $rest = $val.Substring(8)
        $cell.Value = "$rest, System"
    }
}
